# "cleaned defensive actions data"
# Houston Away defensive actions sheet: collapse the pandas-style two-row
# multi-index header into a single visible header row, hide the old
# (now redundant) second header row plus the blank spacer row and the
# "16 Players" totals row, un-merge the old grouped header cells, fill in
# a few zero values that pandas had left blank, and leave the selection
# where the author left it (O21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Undo the old header grouping merges (Tackles/Challenges/Blocks
#     spanned groups) FIRST so every column in row 1 becomes independently
#     writable (a merged range only accepts writes on its top-left cell) ---
[void]$ws.Range("H1:L1").UnMerge()
[void]$ws.Range("M1:P1").UnMerge()
[void]$ws.Range("Q1:S1").UnMerge()

# --- Row 1: replace the old multi-index / merged header labels with the
#     real single-row header (mirrors what used to live in row 2) ---
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2 keeps the old per-column header text (now redundant with row
#     1) but is hidden rather than deleted ---
$ws.Rows(2).Hidden = $true

# --- Row 3 was an empty spacer row; keep it, hidden ---
$ws.Rows(3).Hidden = $true

# --- Fill in the blank Tkl% cells (blank when Tkl was 0) with explicit
#     zeros so every data row is complete ---
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("O19").Value = 0

# --- The "16 Players" totals row is summary data, no longer shown ---
$ws.Rows(20).Hidden = $true

# --- Restore the author's last selection ---
[void]$ws.Range("O21").Select()
